# Weekly driver report update for 2025-04-28
# - Removes the retired "23.40.0.4" bad-driver row (entire row 4), which
#   shifts the Totals row and the whole "Good Drivers" table up by one row.
# - Refreshes the Critical Minutes / Good Roaming % for the remaining bad
#   driver (23.60.0.10) and the Totals row underneath it.
# - Refreshes Total Samples counts for several drivers in the "Good Drivers"
#   table (names/order/dates unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the "23.40.0.4" bad-driver entry entirely (row 4). This shifts
#    every row below it up by one, matching the diff's row renumbering.
$ws.Rows("4:4").Delete()

# 2) Update the remaining bad driver row (23.60.0.10, row 3): Critical
#    Minutes and Good Roaming Calculation (%).
$ws.Range("C3").Value = 2332
$ws.Range("D3").Value = 93.59999999999999

# 3) Update the Totals row (now row 4 after the deletion above).
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = 2332

# 4) Update Total Samples for several "Good Drivers" entries (rows shifted
#    up by one vs. the original file, since row 4 was removed).
$ws.Range("B14").Value = 338880   # 22.230.0.8
$ws.Range("B15").Value = 143869   # 22.200.0.6
$ws.Range("B19").Value = 68450    # 22.10.0.7
$ws.Range("B22").Value = 90508    # 21.40.2.2
$ws.Range("B24").Value = 52515    # 21.10.1.2
